$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as TEXT (mirrors typing into a cell
# that's formatted as Text, or prefixing the entry with an apostrophe),
# then restores the 'Normal' style so no residual number format lingers
# on the cell itself.
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '61.068.03'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '2.400.84'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.49%  '
Set-TextValue 'D5' '566.70'
$ws.Range('E5').Value = '  -0.53%  '
Set-TextValue 'D6' '141.95'
$ws.Range('E6').Value = '  +1.79%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  +2.20%  '
$ws.Range('D9').Value = '2.408.43'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('E11').Value = '  -0.28%  '
Set-TextValue 'D12' '5.19'
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('E13').Value = '  +2.64%  '
Set-TextValue 'D14' '26.48'
$ws.Range('E14').Value = '  +1.47%  '
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '60.810.51'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').Value = '2.423.34'
$ws.Range('E18').Value = '  +0.74%  '
Set-TextValue 'D19' '8.05'
$ws.Range('E19').Value = '  +2.46%  '
Set-TextValue 'D20' '10.67'
$ws.Range('E20').Value = '  +0.51%  '
Set-TextValue 'D21' '324.30'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('E22').Value = '  +1.11%  '
Set-TextValue 'D23' '6.09'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  -0.28%  '
Set-TextValue 'D25' '1.91'
$ws.Range('E25').Value = '  +4.76%  '
Set-TextValue 'D26' '65.21'
$ws.Range('E26').Value = '  +0.74%  '
Set-TextValue 'D27' '587.20'
$ws.Range('E27').Value = '  +0.81%  '
Set-TextValue 'D28' '8.23'
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').Value = '0.0₃0947'
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('E30').Value = '  -0.76%  '
Set-TextValue 'D31' '8.02'
$ws.Range('E31').Value = '  +2.27%  '
$ws.Range('E32').Value = '  +1.56%  '
$ws.Range('E33').Value = '  -0.81%  '
Set-TextValue 'D34' '0.133'
$ws.Range('E34').Value = '  +0.54%  '
Set-TextValue 'D35' '1.48'
$ws.Range('E35').Value = '  +4.12%  '
$ws.Range('E36').Value = '  -0.47%  '
Set-TextValue 'D37' '153.28'
$ws.Range('E37').Value = '  +0.94%  '
Set-TextValue 'D38' '0.372'
$ws.Range('E38').Value = '  +1.06%  '
Set-TextValue 'D39' '4.62'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  +1.05%  '
Set-TextValue 'D44' '41.86'
$ws.Range('E44').Value = '  +1.86%  '
$ws.Range('E45').Value = '  +6.30%  '
$ws.Range('D46').Value = '0.0₆0281'
$ws.Range('E46').Value = '  +3.20%  '
Set-TextValue 'D47' '141.44'
$ws.Range('E47').Value = '  -0.96%  '
Set-TextValue 'D48' '3.52'
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('E49').Value = '  +0.78%  '
$ws.Range('B50').Value = 'Hedera'
$ws.Range('C50').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D50' '0.0510'
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D51' '19.65'
$ws.Range('E51').Value = '  +0.74%  '
